# This workbook holds weekly price-report rows (2..31). The update re-pairs
# the "report" columns (Fecha, Volumen, Precio minimo, Precio maximo,
# Precio promedio ponderado, Origen, Precio $/Kg) onto different rows while
# leaving the identifying columns (Mercado, Region, Codreg, Categoria,
# Variedad, Calidad, Unidad de comercializacion, Kg o Unidades,
# Clasificacion) untouched. In effect, for each destination row we take
# those report columns from a specific source row (a permutation of
# rows 2..31).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# destinationRow = sourceRow (value taken from sourceRow, BEFORE any edits,
# and written into destinationRow)
$rowMap = @{
    2  = 30
    3  = 2
    4  = 19
    5  = 31
    6  = 16
    7  = 27
    8  = 24
    9  = 17
    10 = 18
    11 = 28
    12 = 8
    13 = 26
    14 = 10
    15 = 20
    16 = 4
    17 = 9
    18 = 23
    19 = 11
    20 = 13
    21 = 25
    22 = 3
    23 = 15
    24 = 12
    25 = 7
    26 = 5
    27 = 6
    28 = 22
    29 = 29
    30 = 14
    31 = 21
}

# Columns that move together as a group per row.
$cols = @("D", "J", "K", "L", "M", "O", "P")

# First, snapshot the "before" values for every row/column we will need,
# since several rows feed each other (permutation cycles).
$snapshot = @{}
foreach ($r in $rowMap.Keys) {
    $src = $rowMap[$r]
    if (-not $snapshot.ContainsKey($src)) {
        $rowVals = @{}
        foreach ($c in $cols) {
            $rowVals[$c] = $ws.Range("$c$src").Value()
        }
        $snapshot[$src] = $rowVals
    }
}

# Now write the snapshotted source values into each destination row.
foreach ($r in $rowMap.Keys) {
    $src = $rowMap[$r]
    $rowVals = $snapshot[$src]
    foreach ($c in $cols) {
        $ws.Range("$c$r").Value = $rowVals[$c]
    }
}
